# edit.ps1
# Applies the resume content edits described by the commit "address change at resume":
#   1. Name line: "ZHEYI(jeffrey) TONG" -> "ZHEYI (jeffrey) TONG"
#   2. Street address: "10800 W 133RD TERRACE" -> "11731 W 118th TERRACE"
#   3. Apt/suite line: "APT 12" -> "831"
#   4. Zip code: "66213" -> "66210"
#   5. Education date range: ", Computer Science, Dec 2010" -> ", Computer Science, Aug 2008 - Dec 2010"
#
# All other differences in the source diff are purely cosmetic run-splits
# (Word inserting spell/grammar proofing-error markers around the same
# text) with no net change to the document's text, so they require no
# action here.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "WARNING: could not find '$find'"
    }
    return $ok
}

Replace-Text "ZHEYI(jeffrey) TONG" "ZHEYI (jeffrey) TONG"
Replace-Text "10800 W 133RD TERRACE" "11731 W 118th TERRACE"
Replace-Text "APT 12" "831"
Replace-Text "66213" "66210"
Replace-Text ", Computer Science, Dec 2010" ", Computer Science, Aug 2008 - Dec 2010"

# The source diff also shows the (invisible, non-content) "_GoBack" bookmark
# -- which Word automatically drops at the location of the most recent edit
# -- having moved from the "Environment: Java, Eclipse, ..." bullet to the
# "...T-SQL, SVN." bullet, landing between "T-SQ" and "L, SVN.". Re-create it
# there to mirror that same last-edit marker; Bookmarks.Add re-uses the name,
# relocating (rather than duplicating) the existing bookmark.
$r = $d.Content
$found = $r.Find.Execute("T-SQL, SVN.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $pos = $r.Start + 4
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
} else {
    Write-Host "WARNING: could not find insertion point for _GoBack bookmark"
}
